# Rename the worksheet from "Sheet1" to "DATA_TEMPLATE" (now sourced from the
# settings file rather than hard-coded) and select the header row A4:C4 so the
# template opens with that row highlighted for the user.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "DATA_TEMPLATE"

$ws.Activate()
$ws.Range("A4:C4").Select()
